# Applies changes described by the commit:
# "cambios en PC Gestion Documental para guardar en R"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update account number in H2
$ws.Range("H2").Value = 9498924883

# Clear the value in T2 (cell keeps its style, just loses its content)
$ws.Range("T2").ClearContents()

# Change A3 value from 3 to 2
$ws.Range("A3").Value = 2

# Update the view: remove frozen/scrolled topLeftCell (reset to A1) and change selection to I6
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("I6").Select()
